$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(41.21407239481599, 50.27997645459339, 98.76064511229146, 98.99635139675415, 98.43184127265854, 97.88301481541673, 97.57482661492014, 96.41570945812209)
    3 = @(37.34131802952838, 50.31986102939189, 99.71014989888644, 99.09653945161313, 98.6960493529439, 98.18202371078335, 97.6380440873115, 96.47750989111242)
    4 = @(40.89248999240458, 50.4669486616876, 98.59284432934837, 98.75992146192954, 98.47107576875848, 98.02337899631405, 97.50420190105956, 96.20085015057667)
    5 = @(41.91689552565555, 50.04879908354403, 98.57567300524754, 98.86447131491084, 98.41837211145697, 97.96005616704474, 97.5778095837375, 96.27002368171118)
    6 = @(42.62867068634834, 50.3738104481995, 98.7917614730486, 98.90412690575128, 98.39599034927427, 98.09453442584716, 97.52955856155765, 96.09100202511857)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # Column B is 2
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}
